$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update weeks-of-work figures for the employee cost table based on prototype results
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 10
$ws.Range("C6").Value = 5

# Rename the role in row 6 to reflect the updated project scope
$ws.Range("A6").Value = "Web Designer for Web Application"

# Move the active selection to reflect where the editor left off
$ws.Range("I7").Select()
